# The cover image at the top of the syllabus table used to be wrapped in a
# hyperlink (pointing at the Vimeo trailer). That wrapper hyperlink is
# removed, leaving the picture as a plain, non-linked inline drawing.
# (Every bookmark id / relationship id that was allocated after that
# hyperlink shifts down by one as a natural side effect of the id no
# longer being consumed -- that's not something to replicate by hand.)

$d = $word.ActiveDocument

# The picture lives alone in the first cell of the first (and only) row
# of the first table -- scope to that cell so we unambiguously grab the
# single hyperlink wrapping the drawing (rather than indexing into the
# document-wide Hyperlinks collection).
$tbl = $d.Tables(1)
$cell = $tbl.Cell(1, 1)

if ($cell.Range.Hyperlinks.Count -gt 0) {
    $pictureLink = $cell.Range.Hyperlinks(1)
    # Clearing Address unlinks the hyperlink field, leaving its contents
    # (the <w:r><w:drawing>...</w:drawing></w:r>) in place without the
    # surrounding <w:hyperlink> element.
    $pictureLink.Address = ""
}
